$d = $word.ActiveDocument

# 1. Activation date change
$d.Content.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2022", 2)

# 2. Empty italic paragraph after "Objetivos" paragraph gets English text added
$d.Paragraphs.Item(7).Range.InsertAfter("Provide students with the knowledge of cell biology necessary to understand the other subjects of the course and the training of the Environmental Engineer.")

# 3. "Programa resumido" Portuguese paragraph - remove leading clause
$d.Content.Find.Execute("Origem e evolução da célula; análise estrutural das células ao microscópio", $true, $false, $false, $false, $false, $true, 1, $false, "Análise estrutural das células ao microscópio", 2)

# 4. "Programa resumido" English paragraph - remove leading clause
$d.Content.Find.Execute("The origin and evolution of the cell; organic molecules", $true, $false, $false, $false, $false, $true, 1, $false, "Organic molecules", 2)

# 5. "Programa" Portuguese paragraph - replace opening bullet text
$d.Content.Find.Execute("- Origem e evolução da célula: conceitos básicos de sistemática e filogenia molecular; características dos três diferentes domínios da vida- Análise", $true, $false, $false, $false, $false, $true, 1, $false, "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise", 2)

# 6. "Programa" English paragraph - replace opening text
$d.Content.Find.Execute("The origins and evolution of the cell: basic concepts of systematic and molecular phylogenetics; characteristics of the three domains of life. Microscope", $true, $false, $false, $false, $false, $true, 1, $false, "Cell structure and evolutionary history: prokaryotic microorganisms andeukaryotic and their evolutionary relationships between the Bacteria, Archaea andEukarya.Microscope", 2)

Write-Host "done"
